$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (Health Sciences stays in place)
$ws.Range("B2").Value = 1.4570599903207599
$ws.Range("C2").Value = 2.9110656149999001
$ws.Range("E2").Value = 1.151

# Row 3 label swaps from Multidisciplinary to Physical Sciences, plus new values
$ws.Range("A3").Value = "Physical Sciences"
$ws.Range("B3").Value = 1.0732656119394199
$ws.Range("C3").Value = 0.76643571866648497
$ws.Range("E3").Value = 0.96399999999999997
$ws.Range("F3").Value = 30.405000000000001

# Row 4 label swaps from Physical Sciences to Multidisciplinary, plus new values
$ws.Range("A4").Value = "Multidisciplinary"
$ws.Range("B4").Value = 2.9746764976890501
$ws.Range("C4").Value = 3.53124288650984
$ws.Range("E4").Value = 1.3680000000000001
$ws.Range("F4").Value = 11.342000000000001

# Row 5 values (Social Sciences stays in place)
$ws.Range("B5").Value = 1.32891674463033
$ws.Range("C5").Value = 1.0383018565402899
$ws.Range("E5").Value = 1.157

# Row 6 values (Life Sciences stays in place)
$ws.Range("B6").Value = 1.1875304845872201
$ws.Range("C6").Value = 0.86201307969554597
$ws.Range("E6").Value = 1.0529999999999999

# Update selection from C9 to C10
$ws.Range("C10").Select()
